$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.388.27"
$ws.Range("E2").Value = "  -1.29%  "
$ws.Range("D3").Value = "2.349.69"
$ws.Range("E3").Value = "  +3.40%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'232.07"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("E6").Value = "  +1.87%  "
$ws.Range("D7").Value = "'65.76"
$ws.Range("E7").Value = "  +3.39%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.84%  "
$ws.Range("D10").Value = "'0.0960"
$ws.Range("E10").Value = "  -3.47%  "
$ws.Range("D11").Value = "'57.00"
$ws.Range("E11").Value = "  -0.87%  "
$ws.Range("D12").Value = "'26.73"
$ws.Range("E12").Value = "  -2.39%  "
$ws.Range("D13").Value = "2.699.87"
$ws.Range("E13").Value = "  +3.30%  "
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "'15.41"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("E16").Value = "  +2.09%  "
$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = "  +1.24%  "
$ws.Range("D18").Value = "2.357.35"
$ws.Range("E18").Value = "  +4.00%  "
$ws.Range("D19").Value = "43.386.34"
$ws.Range("E19").Value = "  -1.18%  "
$ws.Range("E20").Value = "  -2.14%  "
$ws.Range("D21").Value = "'74.10"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("E22").Value = "  +1.73%  "
$ws.Range("D23").Value = "'249.42"
$ws.Range("E23").Value = "  -1.26%  "
$ws.Range("D24").Value = "'3.86"
$ws.Range("E24").Value = "  +16.81%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -1.12%  "
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "'9.92"
$ws.Range("E28").Value = "  -1.89%  "
$ws.Range("D29").Value = "'174.98"
$ws.Range("D30").Value = "'22.20"
$ws.Range("E30").Value = "  +6.10%  "
$ws.Range("E31").Value = "  +7.07%  "
$ws.Range("E32").Value = "  -7.26%  "
$ws.Range("D33").Value = "'0.125"
$ws.Range("E33").Value = "  +0.55%  "
$ws.Range("E34").Value = "  +3.83%  "
$ws.Range("D35").Value = "'0.0687"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").Value = "'2.54"
$ws.Range("E37").Value = "  +9.41%  "
$ws.Range("D38").Value = "'6.45"
$ws.Range("E38").Value = "  -1.26%  "
$ws.Range("D39").Value = "'3.62"
$ws.Range("E39").Value = "  -4.99%  "
$ws.Range("E40").Value = "  -3.27%  "
$ws.Range("D41").Value = "'9.02"
$ws.Range("E41").Value = "  +8.97%  "
$ws.Range("E42").Value = "  +0.08%  "
$ws.Range("D43").Value = "'18.07"
$ws.Range("E43").Value = "  +2.42%  "
$ws.Range("E44").Value = "  +9.06%  "
$ws.Range("D45").Value = "'99.08"
$ws.Range("E45").Value = "  +0.73%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'4.39"
$ws.Range("E47").Value = "  +0.45%  "
$ws.Range("E48").Value = "  -4.04%  "
$ws.Range("D49").Value = "1.439.99"
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.574.33"
$ws.Range("E50").Value = "  +3.56%  "
$ws.Range("B51").Value = "TerraClassic"
$ws.Range("C51").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D51").Value = "'0.000204"
$ws.Range("E51").Value = "  -10.44%  "
